$d = $word.ActiveDocument

# 1. Title: merge "Relatório" + " / Comunicado" + " de Não Conformidades"
#    into a single run reading "Relatório de Não Conformidades" by
#    deleting the " / Comunicado" fragment.
$d.Content.Find.Execute(" / Comunicado", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 2)

# 2. Drop the trailing space inside the corrective-action placeholder text.
$d.Content.Find.Execute("conformidade. > ", $false, $false, $false, $false, $false, `
    $true, 1, $false, "conformidade. >", 2)

# 3. Mark the "Normal Table" (Tabelanormal) style as a Quick Style
#    (adds <w:qFormat/> to its style definition).
$d.Styles("Tabelanormal").QuickStyle = $true
